$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# --- New poll data pushed into the rolling window (B6:G8) ---
# New "Latest Morgan" poll added at row 6; previous rows 6 and 7 shift down
# to rows 7 ("Second Morgan") and 8 ("Third Morgan") respectively. The old
# row 8 ("Third Morgan") data rolls off the window.
$row6 = @(57, 59, 57.5, 51.5, 53.5, 59.5)
$row7 = @(56.5, 54, 59, 51.5, 55.5, 64)
$row8 = @(56, 58, 59, 48.5, 50.5, 60.5)

# --- Updated Essential (3 avg.) rows ---
$row12 = @(52.142857142857139, 49.635036496350367, 52.857142857142861, 51.957295373665481, 54.255319148936174, 50)
$row13 = @(51.624548736462089, 47.653429602888082, 50.719424460431661, 53.429602888086642, 59.574468085106382, 50.537634408602152)

$cols = @("B", "C", "D", "E", "F", "G")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])6").Value = $row6[$i]
    $ws.Range("$($cols[$i])7").Value = $row7[$i]
    $ws.Range("$($cols[$i])8").Value = $row8[$i]
    $ws.Range("$($cols[$i])12").Value = $row12[$i]
    $ws.Range("$($cols[$i])13").Value = $row13[$i]
}

# --- Re-enter the swing-deviation formulas as filled ranges so Excel stores ---
# --- them as shared formulas (matches a fill-down/fill-right re-entry) ---
$ws.Range("C16:G19").Formula = "=(C5-C`$2)-(`$B5-`$B`$2)"
$ws.Range("C20:E24").Formula = "=(C9-C`$2)-(`$B9-`$B`$2)"

# --- Update the active selection on the Calc sheet ---
$ws.Range("K8").Select()

$wb.Save()
